# recentrifuge/test/mock.xlsx
# "Add native and contaminated strains to mock for challenging ROC"
#
# - Update the Homo sapiens (row 2) counts.
# - Insert two new taxa rows (Methanobacterium formicicum DSM 3637 / JCM 10132)
#   right before "Methylobacterium radiotolerans", pushing the rest of the
#   table (and the trailing SUM row) down by two rows.
# - Leave the active selection on the second new row, as in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mock")

# --- Row 2 (Homo sapiens): updated sample counts ---------------------------
$ws.Range("C2").Value = 86830
$ws.Range("D2").Value = 60573
$ws.Range("E2").Value = 87100
$ws.Range("F2").Value = 31168
$ws.Range("G2").Value = 27994
$ws.Range("H2").Value = 17940
$ws.Range("I2").Value = 17393

# --- Insert two fresh rows before "Methylobacterium radiotolerans" (row 19) ---
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(19).Insert()

# New rows inherit the formatting of the row above on insert; restore the
# same fills used elsewhere in the sheet for this kind of entry (row 13 style
# for row 19, row 8 style for row 20) via a format-only paste.
$ws.Range("A13").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A8").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 19: Methanobacterium formicicum DSM 3637
$ws.Range("A19").Value = "Methanobacterium formicicum DSM 3637"
$ws.Range("B19").Value = 1204725
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 3
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 300
$ws.Range("I19").Value = 3000

# Row 20: Methanobacterium formicicum JCM 10132
$ws.Range("A20").Value = "Methanobacterium formicicum JCM 10132"
$ws.Range("B20").Value = 1300163
$ws.Range("C20").Value = 300
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 30
$ws.Range("F20").Value = 3000
$ws.Range("G20").Value = 300
$ws.Range("H20").Value = 30
$ws.Range("I20").Value = 3

# --- Column A is now a little wider to fit the longer species names --------
$ws.Columns.Item(1).ColumnWidth = 34.998697916666664

# --- Match the saved selection (active cell on the new JCM 10132 row) ------
[void]$ws.Range("B20").Select()

Write-Output "edit applied"
